$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.960.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.867.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.73%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.20%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5089"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3938"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08200"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.36%  "
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.093"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.874.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.282"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.178"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.02%  "
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001087"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06391"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.952.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.822"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.08%  "
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.171"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.079.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.223"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.924"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.730"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02427"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.208"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06344"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2139"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.172"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.492"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6300"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.05%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.200"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5895"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.003"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.203"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.117"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.14%  "
